$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 46; this shifts existing rows 46-106 down to 47-107
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new data entry
$ws.Cells.Item(46, 1).Value = 8
$ws.Cells.Item(46, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(46, 3).Value = "Coquimbo"
$ws.Cells.Item(46, 4).Value = 44771
$ws.Cells.Item(46, 5).Value = 4
$ws.Cells.Item(46, 6).Value = 100112052
$ws.Cells.Item(46, 7).Value = "Albahaca"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 1000
$ws.Cells.Item(46, 11).Value = 3500
$ws.Cells.Item(46, 12).Value = 4000
$ws.Cells.Item(46, 13).Value = 3750
$ws.Cells.Item(46, 14).Value = "`$/paquete"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 3750
$ws.Cells.Item(46, 17).Value = 1
$ws.Cells.Item(46, 18).Value = "Hortaliza"
